# Update NATMI ligand-receptor (Nodal -> Acvr2b) sheet with recomputed
# TPM-derived expression/specificity values for rows 2-17 (per the upstream
# scripts' new TPM matrix). Only numeric value columns E:T change; the
# cluster-label columns A:D are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => @{ column letter = new value }, taken from the updated TPM output
$newValues = @{
    2 = @{ "E" = 3; "F" = 1; "G" = 1.591135333333333; "H" = 4.773406; "I" = 0.4707829234247397; "J" = 0.4707829234247397; "K" = 3; "L" = 1; "M" = 2.110264333333333; "N" = 6.330793; "O" = 0.3832041185227171; "P" = 0.3832041185227171; "Q" = 3.357716143439777; "R" = 30.219445290958; "S" = 0.1804059551865252; "T" = 0.1804059551865252 }
    3 = @{ "E" = 3; "F" = 1; "G" = 1.591135333333333; "H" = 4.773406; "I" = 0.4707829234247397; "J" = 0.4707829234247397; "O" = 0.0946183755984393; "P" = 0.0946183755984393; "Q" = 0.8290663692177775; "R" = 7.461597322959999; "S" = 0.04454471547393331; "T" = 0.04454471547393331 }
    4 = @{ "E" = 3; "F" = 1; "G" = 1.591135333333333; "H" = 4.773406; "I" = 0.4707829234247397; "J" = 0.4707829234247397; "M" = 2.065388333333333; "N" = 6.196165; "O" = 0.3750550597762889; "P" = 0.3750550597762889; "Q" = 3.286312354221111; "R" = 29.57681118799; "S" = 0.1765695174867218; "T" = 0.1765695174867218 }
    5 = @{ "E" = 3; "F" = 1; "G" = 1.591135333333333; "H" = 4.773406; "I" = 0.4707829234247397; "J" = 0.4707829234247397; "M" = 0.8101876666666666; "N" = 2.430563; "O" = 0.1471224461025547; "P" = 0.1471224461025547; "Q" = 1.289118223064222; "R" = 11.602064007578; "S" = 0.06926273527755941; "T" = 0.06926273527755943 }
    6 = @{ "I" = 0.3035973020998604; "J" = 0.3035973020998604; "K" = 3; "L" = 1; "M" = 2.110264333333333; "N" = 6.330793; "O" = 0.3832041185227171; "P" = 0.3832041185227171; "Q" = 2.165315502418444; "R" = 19.487839521766; "S" = 0.1163397365370521; "T" = 0.1163397365370521 }
    7 = @{ "I" = 0.3035973020998604; "J" = 0.3035973020998604; "O" = 0.0946183755984393; "P" = 0.0946183755984393; "S" = 0.02872588356075744; "T" = 0.02872588356075744 }
    8 = @{ "I" = 0.3035973020998604; "J" = 0.3035973020998604; "M" = 2.065388333333333; "N" = 6.196165; "O" = 0.3750550597762889; "P" = 0.3750550597762889; "Q" = 2.119268807247778; "R" = 19.07341926523; "S" = 0.1138657042869832; "T" = 0.1138657042869832 }
    9 = @{ "I" = 0.3035973020998604; "J" = 0.3035973020998604; "M" = 0.8101876666666666; "N" = 2.430563; "O" = 0.1471224461025547; "P" = 0.1471224461025547; "Q" = 0.8313233023895555; "R" = 7.481909721506; "S" = 0.04466597771506774; "T" = 0.04466597771506774 }
    10 = @{ "G" = 0.730693; "H" = 2.192079; "I" = 0.2161964349979826; "J" = 0.2161964349979826; "K" = 3; "L" = 1; "M" = 2.110264333333333; "N" = 6.330793; "O" = 0.3832041185227171; "P" = 0.3832041185227171; "Q" = 1.541955376516333; "R" = 13.877598388647; "S" = 0.08284736430115583; "T" = 0.08284736430115583 }
    11 = @{ "G" = 0.730693; "H" = 2.192079; "I" = 0.2161964349979826; "J" = 0.2161964349979826; "O" = 0.0946183755984393; "P" = 0.0946183755984393; "Q" = 0.3807300232933333; "R" = 3.42657020964; "S" = 0.02045615548968269; "T" = 0.02045615548968269 }
    12 = @{ "G" = 0.730693; "H" = 2.192079; "I" = 0.2161964349979826; "J" = 0.2161964349979826; "M" = 2.065388333333333; "N" = 6.196165; "O" = 0.3750550597762889; "P" = 0.3750550597762889; "Q" = 1.509164797448334; "R" = 13.582483177035; "S" = 0.08108556685158892; "T" = 0.08108556685158892 }
    13 = @{ "G" = 0.730693; "H" = 2.192079; "I" = 0.2161964349979826; "J" = 0.2161964349979826; "M" = 0.8101876666666666; "N" = 2.430563; "O" = 0.1471224461025547; "P" = 0.1471224461025547; "Q" = 0.5919984567196667; "R" = 5.327986110477; "S" = 0.03180734835555517; "T" = 0.03180734835555517 }
    14 = @{ "G" = 0.03184866666666667; "H" = 0.09554600000000001; "I" = 0.009423339477417213; "J" = 0.009423339477417213; "K" = 3; "L" = 1; "M" = 2.110264333333333; "N" = 6.330793; "O" = 0.3832041185227171; "P" = 0.3832041185227171; "Q" = 0.06720910533088889; "R" = 0.604881947978; "S" = 0.003611062497983984; "T" = 0.003611062497983984 }
    15 = @{ "G" = 0.03184866666666667; "H" = 0.09554600000000001; "I" = 0.009423339477417213; "J" = 0.009423339477417213; "O" = 0.0946183755984393; "P" = 0.0946183755984393; "Q" = 0.01659485392888889; "R" = 0.14935368536; "S" = 0.0008916210740658626; "T" = 0.0008916210740658626 }
    16 = @{ "G" = 0.03184866666666667; "H" = 0.09554600000000001; "I" = 0.009423339477417213; "J" = 0.009423339477417213; "M" = 2.065388333333333; "N" = 6.196165; "O" = 0.3750550597762889; "P" = 0.3750550597762889; "Q" = 0.06577986456555557; "R" = 0.5920187810900001; "S" = 0.003534271150994975; "T" = 0.003534271150994975 }
    17 = @{ "G" = 0.03184866666666667; "H" = 0.09554600000000001; "I" = 0.009423339477417213; "J" = 0.009423339477417213; "M" = 0.8101876666666666; "N" = 2.430563; "O" = 0.1471224461025547; "P" = 0.1471224461025547; "Q" = 0.02580339693311112; "R" = 0.232230572398; "S" = 0.00138638475437239; "T" = 0.00138638475437239 }
}

foreach ($row in $newValues.Keys) {
    foreach ($col in $newValues[$row].Keys) {
        $ws.Range("$col$row").Value = $newValues[$row][$col]
    }
}
